$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 34: new data row (previously a blank row formatted D34/E34/F34) ---
$ws.Range("A34").Value = 2014
$ws.Range("B34").Value = 3
$ws.Range("C34").Value = 2
$ws.Range("D34").Value = 0.40277777777777773
$ws.Range("E34").Value = 0.40972222222222227
$ws.Range("F34").Formula = "=(E34-D34)*24*60"
$ws.Range("G34").Formula = "=F34/60"

# --- Row 35: new data row (previously held the "sum [min]" label) ---
$ws.Range("E35").Style = "Normal"
$ws.Range("A35").Value = 2014
$ws.Range("B35").Value = 3
$ws.Range("C35").Value = 2
$ws.Range("D35").Value = 0.4201388888888889
$ws.Range("E35").Value = 0.5
$ws.Range("F35").Formula = "=(E35-D35)*24*60"
$ws.Range("G35").Formula = "=F35/60"

# --- Row 36: new blank formatted row (mirrors previous blank row 34) ---
$ws.Range("E36").Style = "Normal"
$ws.Range("F36").Style = "Normal"
$ws.Range("D36").Value = ""
$ws.Range("E36").Value = ""
$ws.Range("F36").Value = ""

# --- Row 37: "sum [min]" label row (content previously on row 35) ---
$ws.Range("E37").Value = "sum [min]"
$ws.Range("F37").Style = "Normal"
$ws.Range("F37").Formula = "=SUM(F2:F36)"

# --- Row 38: "sum [h]" label row (content previously on row 36) ---
$ws.Range("E38").Value = "sum [h]"
$ws.Range("F38").Formula = "=F37/60"

# --- Row 39: "sum [working weeks]" label row (content previously on row 37) ---
$ws.Range("E39").Value = "sum [working weeks]"
$ws.Range("F39").Formula = "=F38/38.5"

# --- Number formats / alignment to match the sheet's existing style classes ---
$ws.Range("D34:E35").NumberFormat = "hh:mm;@"
$ws.Range("F34:F35").NumberFormat = "0"
$ws.Range("G34:G35").NumberFormat = "0.00"

$ws.Range("D36:E36").NumberFormat = "hh:mm;@"
$ws.Range("F36").NumberFormat = "0"

$ws.Range("F37").NumberFormat = "0"
$ws.Range("F38:F39").NumberFormat = "0.00"

$ws.Range("E37:E39").HorizontalAlignment = -4152

$ws.Range("F35").Select()
